$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.112.79'
$ws.Range('E2').Value = '  -1.50%  '

$ws.Range('D3').Value = '1.656.20'
$ws.Range('E3').Value = '  -1.38%  '

$cell = $ws.Range('D5')
$cell.Value = "'" + '216.12'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -1.74%  '

$cell = $ws.Range('D6')
$cell.Value = "'" + '0.5167'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  -2.79%  '

$ws.Range('E7').Value = '  +0.36%  '

$cell = $ws.Range('D8')
$cell.Value = "'" + '0.2625'
$cell.Style = "Normal"
$ws.Range('E8').Value = '  -2.90%  '

$cell = $ws.Range('D9')
$cell.Value = "'" + '0.06257'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  -2.48%  '

$cell = $ws.Range('D10')
$cell.Value = "'" + '20.71'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -5.19%  '

$ws.Range('E11').Value = '  -1.29%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.657.63'
$ws.Range('E12').Value = '  -1.46%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range('D13')
$cell.Value = "'" + '4.420'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  -2.12%  '

$ws.Range('D14').Value = '1.882.77'
$ws.Range('E14').Value = '  -1.37%  '

$cell = $ws.Range('D15')
$cell.Value = "'" + '0.5404'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -3.49%  '

$ws.Range('D16').Value = '0.0₅8104'
$ws.Range('E16').Value = '  -3.03%  '

$cell = $ws.Range('D17')
$cell.Value = "'" + '64.70'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  -1.67%  '

$ws.Range('D18').Value = '26.146.31'
$ws.Range('E18').Value = '  -1.48%  '

$ws.Range('E19').Value = '  +0.38%  '

$cell = $ws.Range('D20')
$cell.Value = "'" + '4.610'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  -4.04%  '

$cell = $ws.Range('D21')
$cell.Value = "'" + '191.24'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  -1.14%  '

$cell = $ws.Range('D22')
$cell.Value = "'" + '10.06'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -2.66%  '

$cell = $ws.Range('D23')
$cell.Value = "'" + '6.017'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  -5.10%  '

$ws.Range('E24').Value = '  +0.41%  '

$cell = $ws.Range('D25')
$cell.Value = "'" + '139.78'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +0.21%  '

$cell = $ws.Range('D26')
$cell.Value = "'" + '0.1223'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  -4.23%  '

$cell = $ws.Range('D27')
$cell.Value = "'" + '7.161'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  -3.55%  '

$cell = $ws.Range('D28')
$cell.Value = "'" + '16.06'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -1.49%  '

$cell = $ws.Range('D29')
$cell.Value = "'" + '1.401'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -2.99%  '

$cell = $ws.Range('D30')
$cell.Value = "'" + '0.05955'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  -5.36%  '

$ws.Range('E31').Value = '  -1.64%  '

$cell = $ws.Range('D32')
$cell.Value = "'" + '3.535'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  -2.22%  '

$cell = $ws.Range('D33')
$cell.Value = "'" + '3.251'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  -6.18%  '

$cell = $ws.Range('D34')
$cell.Value = "'" + '1.602'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -5.68%  '

$cell = $ws.Range('D35')
$cell.Value = "'" + '0.9635'
$cell.Style = "Normal"
$ws.Range('E35').Value = '  -5.05%  '

$cell = $ws.Range('D36')
$cell.Value = "'" + '2.427'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  +0.25%  '

$cell = $ws.Range('D37')
$cell.Value = "'" + '2.774'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  -0.58%  '

$cell = $ws.Range('D38')
$cell.Value = "'" + '0.5663'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  -8.70%  '

$cell = $ws.Range('D39')
$cell.Value = "'" + '0.01590'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  -2.81%  '

$cell = $ws.Range('D40')
$cell.Value = "'" + '5.962'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  -3.13%  '

$cell = $ws.Range('D41')
$cell.Value = "'" + '0.8545'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -0.98%  '

$ws.Range('E42').Value = '  +0.38%  '

$ws.Range('D43').Value = '1.008.70'
$ws.Range('E43').Value = '  -8.04%  '

$cell = $ws.Range('D44')
$cell.Value = "'" + '100.38'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  -0.35%  '

$ws.Range('D45').Value = '1.797.69'
$ws.Range('E45').Value = '  -1.43%  '

$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  -1.56%  '

$cell = $ws.Range('D47')
$cell.Value = "'" + '56.55'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -3.87%  '

$cell = $ws.Range('D48')
$cell.Value = "'" + '1.005'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -0.02%  '

$ws.Range('E49').Value = '  -2.52%  '

$cell = $ws.Range('D50')
$cell.Value = "'" + '0.05173'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  -0.43%  '

$cell = $ws.Range('D51')
$cell.Value = "'" + '0.4198'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  -0.85%  '

